$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: sender/recipient data (keep A3,B3; change C3,D3,F3) ---
$ws.Cells.Item(3,3).Value = "314b065159e8e9c"
$ws.Cells.Item(3,4).Value = "RS35155000000000563774"
$ws.Cells.Item(3,6).Value = "EUR"

# --- Row 4: clear A4, change B4-D4, keep E4/G4, change F4, H4 ---
$ws.Cells.Item(4,1).ClearContents()
$ws.Cells.Item(4,2).Value = "BE39914001921319"
$ws.Cells.Item(4,3).Value = "314b065159e8e9c"
$ws.Cells.Item(4,4).Value = "RS35155000000000563774"
$ws.Cells.Item(4,6).Value = "EUR"
$ws.Cells.Item(4,8).Value = "recipient"

# --- Row 5: keep A5; change B5; clear C5; keep D5 text but change value; change E5; clear F5,G5,H5 ---
$ws.Cells.Item(5,2).Value = "BE39914001921319"
$ws.Cells.Item(5,3).ClearContents()
$ws.Cells.Item(5,4).Value = "RS35155000000000563774"
$ws.Cells.Item(5,5).Value = 11
$ws.Cells.Item(5,6).ClearContents()
$ws.Cells.Item(5,7).ClearContents()
$ws.Cells.Item(5,8).ClearContents()

# --- Row 11: clear A11:J11, keep K11 ---
$ws.Range("A11:J11").ClearContents()

# --- Delete rows 12 and 13 (shift rows up) ---
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(12).Delete()
